# Review until chapter 3
# Update the two SmartArt node labels on slide 1's diagram:
#   "گره‌های اجرایی" -> "طرح اجرایی"
#   "گراف اجرا"      -> "طرح منطقی"
# (The diagram's text lives in both the diagramData and diagramDrawing
# parts; editing through the SmartArt object model keeps them in sync.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$sa = $sh.SmartArt
$nodes = $sa.AllNodes

for ($i = 1; $i -le $nodes.Count; $i++) {
    $node = $nodes.Item($i)
    $t = $node.TextFrame2.TextRange.Text
    if ($t -eq "گره‌های اجرایی") {
        $node.TextFrame2.TextRange.Text = "طرح اجرایی"
    } elseif ($t -eq "گراف اجرا") {
        $node.TextFrame2.TextRange.Text = "طرح منطقی"
    }
}
